{"js": "// The diff appends new sentences to the end of the second paragraph\n// (\"In my opinion,it's good to do some reading at home on rainy days.\")\n// so that it reads:\n// \"In my opinion,it's good to do some reading at home on rainy days.\n//  I think I don't need to emphasize the importance ofreading as\n//  everybody knows.Look back our daily life,how many people have time\n//  to read now and then?Why don't we catch this chance to kill time and\n//  enhance ourselves?Secondly,\"\n//\n// Locate the end of that paragraph robustly (search for the unique\n// trailing phrase instead of assuming a fixed paragraph index), then\n// append the new text. The closing right single quote (\u2019) that Word\n// inserts for \"don't\" is kept in its own run, matching the convention\n// already used elsewhere in this document (e.g. \"Don\u2019t\", \"it\u2019s\").\n\nconst anchorText = \"on rainy days.\";\nconst results = context.document.body.search(anchorText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the paragraph ending in 'on rainy days.'\");\n}\n\nconst match = results.items[0];\nconst paragraphs = match.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = paragraphs.items[0];\n\nconst rightSingleQuote = \"\\u2019\";\nconst beforeApostrophe1 = \"I think I don\";\nconst betweenApostrophes = \"t need to emphasize the importance ofreading as everybody knows.\" +\n  \"Look back our daily life,how many people have time to read now and then?Why don\";\nconst afterApostrophe2 = \"t we catch this chance to kill time and enhance ourselves?Secondly,\";\n\ntargetParagraph.insertText(beforeApostrophe1, \"End\");\ntargetParagraph.insertText(rightSingleQuote, \"End\");\ntargetParagraph.insertText(betweenApostrophes, \"End\");\ntargetParagraph.insertText(rightSingleQuote, \"End\");\ntargetParagraph.insertText(afterApostrophe2, \"End\");\n\nawait context.sync();\n", "ps1": "# The diff appends new sentences to the end of the second paragraph\n# (\"In my opinion,it's good to do some reading at home on rainy days.\")\n# so that it reads:\n# \"In my opinion,it's good to do some reading at home on rainy days.\n#  I think I don't need to emphasize the importance ofreading as\n#  everybody knows.Look back our daily life,how many people have time\n#  to read now and then?Why don't we catch this chance to kill time and\n#  enhance ourselves?Secondly,\"\n#\n# Locate the end of that paragraph robustly with Find (instead of\n# assuming a fixed paragraph index), then append the new text there.\n# The closing right single quote (U+2019) that Word uses for \"don't\" is\n# inserted as its own piece, matching the convention already used\n# elsewhere in this document (e.g. \"Don't\", \"it's\") where the curly\n# apostrophe sits in its own run.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"on rainy days.\")\nif (-not $found) {\n    throw \"Could not find the paragraph ending in 'on rainy days.'\"\n}\n\n# Collapse the found range to its end (just before the paragraph mark)\n# so the new text is appended right after \"...on rainy days.\"\n$rng.Collapse(0)\n\n$rightSingleQuote = [char]0x2019\n\n$rng.InsertAfter(\"I think I don\")\n$rng.InsertAfter($rightSingleQuote)\n$rng.InsertAfter(\"t need to emphasize the importance ofreading as everybody knows.Look back our daily life,how many people have time to read now and then?Why don\")\n$rng.InsertAfter($rightSingleQuote)\n$rng.InsertAfter(\"t we catch this chance to kill time and enhance ourselves?Secondly,\")\n"}
